{"js": "const replacements = [\n  [\"2025-09-02 Tuesday\", \"2025-09-03 Wednesday\"],\n  [\"51\u00d784=4284\", \"33\u00d747=1551\"],\n  [\"51\u00d745=2295\", \"13\u00d739=507\"],\n  [\"53\u00d790=4770\", \"42\u00d750=2100\"],\n  [\"57\u00d711=627\", \"87\u00d751=4437\"],\n  [\"81\u00d784=6804\", \"59\u00d779=4661\"],\n  [\"26\u00d797=2522\", \"23\u00d729=667\"],\n  [\"96\u00d775=7200\", \"34\u00d724=816\"],\n  [\"24\u00d757=1368\", \"13\u00d742=546\"],\n  [\"84\u00d748=4032\", \"55\u00d765=3575\"],\n  [\"21\u00d735=735\", \"92\u00d711=1012\"],\n  [\"54\u00d715=810\", \"47\u00d714=658\"],\n  [\"39\u00d749=1911\", \"83\u00d755=4565\"],\n  [\"20\u00d765=1300\", \"19\u00d728=532\"],\n  [\"91\u00d797=8827\", \"16\u00d784=1344\"],\n  [\"30\u00d737=1110\", \"95\u00d757=5415\"],\n  [\"89\u00d780=7120\", \"95\u00d735=3325\"],\n  [\"50\u00d712=600\", \"14\u00d757=798\"],\n  [\"54\u00d796=5184\", \"93\u00d722=2046\"],\n  [\"35\u00d751=1785\", \"14\u00d796=1344\"],\n  [\"15\u00d787=1305\", \"39\u00d737=1443\"],\n  [\"17\u00d745=765\", \"34\u00d712=408\"],\n  [\"53\u00d747=2491\", \"23\u00d740=920\"],\n  [\"45\u00d797=4365\", \"67\u00d730=2010\"],\n  [\"27\u00d727=729\", \"28\u00d722=616\"],\n  [\"87\u00d742=3654\", \"65\u00d749=3185\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-09-02 Tuesday\", \"2025-09-03 Wednesday\"),\n    @(\"51\u00d784=4284\", \"33\u00d747=1551\"),\n    @(\"51\u00d745=2295\", \"13\u00d739=507\"),\n    @(\"53\u00d790=4770\", \"42\u00d750=2100\"),\n    @(\"57\u00d711=627\", \"87\u00d751=4437\"),\n    @(\"81\u00d784=6804\", \"59\u00d779=4661\"),\n    @(\"26\u00d797=2522\", \"23\u00d729=667\"),\n    @(\"96\u00d775=7200\", \"34\u00d724=816\"),\n    @(\"24\u00d757=1368\", \"13\u00d742=546\"),\n    @(\"84\u00d748=4032\", \"55\u00d765=3575\"),\n    @(\"21\u00d735=735\", \"92\u00d711=1012\"),\n    @(\"54\u00d715=810\", \"47\u00d714=658\"),\n    @(\"39\u00d749=1911\", \"83\u00d755=4565\"),\n    @(\"20\u00d765=1300\", \"19\u00d728=532\"),\n    @(\"91\u00d797=8827\", \"16\u00d784=1344\"),\n    @(\"30\u00d737=1110\", \"95\u00d757=5415\"),\n    @(\"89\u00d780=7120\", \"95\u00d735=3325\"),\n    @(\"50\u00d712=600\", \"14\u00d757=798\"),\n    @(\"54\u00d796=5184\", \"93\u00d722=2046\"),\n    @(\"35\u00d751=1785\", \"14\u00d796=1344\"),\n    @(\"15\u00d787=1305\", \"39\u00d737=1443\"),\n    @(\"17\u00d745=765\", \"34\u00d712=408\"),\n    @(\"53\u00d747=2491\", \"23\u00d740=920\"),\n    @(\"45\u00d797=4365\", \"67\u00d730=2010\"),\n    @(\"27\u00d727=729\", \"28\u00d722=616\"),\n    @(\"87\u00d742=3654\", \"65\u00d749=3185\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $result) {\n        Write-Output \"FAILED to replace: $old\"\n    }\n}\n"}
